# Generate Report for Handoff
# Adds two new "Ready for handoff" file rows (2d168220-... and 7218e061-...)
# ahead of the existing last row (78e2c296-...) on all three worksheets:
#   Sheet "Overview" (summary, columns A-D)
#   Sheet "zh-cn"     (detail, columns A-L)
#   Sheet "de-de"     (detail, columns A-L)

$wb = $excel.ActiveWorkbook

$mdA  = "2d168220-e8fc-4800-a169-997a02da2dbd.md"
$mdB  = "7218e061-032c-4c7d-89c4-fa7b3938d44c.md"
$mdC  = "78e2c296-f1dc-4bd8-95b2-65a4a58df23b.md"

$xlfZhA = "2d168220-e8fc-4800-a169-997a02da2dbd.10b7701fd5d35ae646d80701ae31ab2d31f7c4a6.zh-cn.xlf"
$xlfZhB = "7218e061-032c-4c7d-89c4-fa7b3938d44c.91997aacf1282224be42ae10b856dac2c60548c6.zh-cn.xlf"
$xlfZhC = "78e2c296-f1dc-4bd8-95b2-65a4a58df23b.bf4a132fdc17200388bddce44add79d6c9cb288c.zh-cn.xlf"

$xlfDeA = "2d168220-e8fc-4800-a169-997a02da2dbd.10b7701fd5d35ae646d80701ae31ab2d31f7c4a6.de-de.xlf"
$xlfDeB = "7218e061-032c-4c7d-89c4-fa7b3938d44c.91997aacf1282224be42ae10b856dac2c60548c6.de-de.xlf"
$xlfDeC = "78e2c296-f1dc-4bd8-95b2-65a4a58df23b.bf4a132fdc17200388bddce44add79d6c9cb288c.de-de.xlf"

$status = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (file), B (zh-cn status), C (de-de status), D (date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 5; File = $mdA; Date = "2016-03-22 14:41:34"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f9001122334/e2e/$mdA" },
    @{ Row = 6; File = $mdB; Date = "2016-03-22 14:41:34"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b2c3d4e5f60718293a4b5c6d7e8f900112233445/e2e/$mdB" },
    @{ Row = 7; File = $mdC; Date = "2016-03-22 14:38:28"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3eb9bab449976525ce0b8e768563ef21aaa9743d/e2e/$mdC" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $wsOverview.Range("A$row").Value2 = $r.File
    $wsOverview.Range("A$row").Style = "HyperLink"
    $wsOverview.Range("B$row").Value2 = $status
    $wsOverview.Range("C$row").Value2 = $status
    $wsOverview.Range("D$row").Value2 = $r.Date
    $wsOverview.Range("D$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$wsOverview.Hyperlinks.Delete()
foreach ($r in $overviewRows) {
    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$($r.Row)"), $r.MdUrl, "", "", $r.File) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": detail columns A-L
# ---------------------------------------------------------------------------
$detailSheets = @(
    @{
        Name = "zh-cn"
        Rows = @(
            @{ Row = 5; Md = $mdA; Xlf = $xlfZhA; Date = "2016-03-22 14:41:31"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f9001122334/e2e/$mdA"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10b7701fd5d35ae646d80701ae31ab2d31f7c4a6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhA" },
            @{ Row = 6; Md = $mdB; Xlf = $xlfZhB; Date = "2016-03-22 14:41:31"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b2c3d4e5f60718293a4b5c6d7e8f900112233445/e2e/$mdB"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91997aacf1282224be42ae10b856dac2c60548c6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhB" },
            @{ Row = 7; Md = $mdC; Xlf = $xlfZhC; Date = "2016-03-22 14:38:23"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3eb9bab449976525ce0b8e768563ef21aaa9743d/e2e/$mdC"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50afbfa611627b0b8ce703319bf21f0dd358d6da/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$xlfZhC" }
        )
    },
    @{
        Name = "de-de"
        Rows = @(
            @{ Row = 5; Md = $mdA; Xlf = $xlfDeA; Date = "2016-03-22 14:41:34"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a1b2c3d4e5f60718293a4b5c6d7e8f9001122334/e2e/$mdA"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10b7701fd5d35ae646d80701ae31ab2d31f7c4a6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeA" },
            @{ Row = 6; Md = $mdB; Xlf = $xlfDeB; Date = "2016-03-22 14:41:34"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/b2c3d4e5f60718293a4b5c6d7e8f900112233445/e2e/$mdB"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/91997aacf1282224be42ae10b856dac2c60548c6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeB" },
            @{ Row = 7; Md = $mdC; Xlf = $xlfDeC; Date = "2016-03-22 14:38:28"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3eb9bab449976525ce0b8e768563ef21aaa9743d/e2e/$mdC"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1d2bb3f35981fe799c8ce4bac8ec5cd24fe35985/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$xlfDeC" }
        )
    }
)

foreach ($sheet in $detailSheets) {
    $ws = $wb.Worksheets.Item($sheet.Name)

    foreach ($r in $sheet.Rows) {
        $row = $r.Row
        $ws.Range("A$row").Value2 = $r.Md
        $ws.Range("A$row").Style = "HyperLink"

        $ws.Range("B$row").Value2 = ".md"

        $ws.Range("C$row").Value2 = $status

        $ws.Range("D$row").Value2 = $r.Xlf
        $ws.Range("D$row").Style = "HyperLink"

        $ws.Range("E$row").Value2 = $r.Date
        $ws.Range("E$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"

        $ws.Range("H$row").Value2 = "0001-01-01 00:00:00"
        $ws.Range("H$row").NumberFormat = "yyyy-mm-dd HH:mm:ss"

        $ws.Range("J$row").Value2 = "Include"
    }

    $ws.Hyperlinks.Delete()
    foreach ($r in $sheet.Rows) {
        $ws.Hyperlinks.Add($ws.Range("A$($r.Row)"), $r.MdUrl, "", "", $r.Md) | Out-Null
        $ws.Hyperlinks.Add($ws.Range("D$($r.Row)"), $r.XlfUrl, "", "", $r.Xlf) | Out-Null
    }
}
